$d = $word.ActiveDocument

# 1) Objetivos paragraph: split after "...ao longo do curso," before "conforme..."
$d.Content.Find.Execute(
    "curso,conforme projeto aprovado",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "curso,^lconforme projeto aprovado", 2)

# 2) Programa resumido paragraph: split after "...constituir-se num" before "projeto de engenharia..."
$d.Content.Find.Execute(
    "constituir-se numprojeto de engenharia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "constituir-se num^lprojeto de engenharia", 2)

# 3) Programa paragraph: split after "(3) a" before "revisão bibliográfica..." and after
#    "resultados, (7)" before "as conclusões..."
$d.Content.Find.Execute(
    "objetivos, (3) arevisão bibliográfica",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "objetivos, (3) a^lrevisão bibliográfica", 2)

$d.Content.Find.Execute(
    "resultados, (7)as conclusões",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "resultados, (7)^las conclusões", 2)
